# Update the "K" column (column G) values on the active worksheet
# with newly calculated s_vals, per commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,0,2,0,1,0,0,1,0,1,2,0,2,1,2,1,0,1,0,3,2,1,0,0,0,2,0,2,1,1,2,1,1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
